$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ceinture_respiration_Main")

# Update the MPN cell for BT1 (row 3) to add the second option part number
$ws.Range("E3").Value = "BH9VPC / 36-232-ND"

# Restore selection to match the saved state
$ws.Range("H18").Select()
